# Scheduled runner update: refresh market-price derived columns (H-N)
# across the Titan_Profits leve-profit tables for ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3475438
$ws.Range("I76").Value = 4276543
$ws.Range("J76").Value = 3983.3333
$ws.Range("K76").Value = 4276543
$ws.Range("L76").Value = 3983.3333
$ws.Range("M76").Value = -4276228
$ws.Range("N76").Value = -4613.3333
$ws.Range("H79").Value = 3475438
$ws.Range("I79").Value = 4276543
$ws.Range("J79").Value = 3983.3333
$ws.Range("K79").Value = 4276543
$ws.Range("L79").Value = 3983.3333
$ws.Range("M79").Value = -4275451
$ws.Range("N79").Value = -6167.3333
$ws.Range("H113").Value = 1987.4166
$ws.Range("I113").Value = 1983.2222
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1983.2222
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1270.7778
$ws.Range("N113").Value = -8508
$ws.Range("H116").Value = 8651600
$ws.Range("I116").Value = 12582247
$ws.Range("J116").Value = 4178
$ws.Range("K116").Value = 12582247
$ws.Range("L116").Value = 4178
$ws.Range("M116").Value = -12578805
$ws.Range("N116").Value = -11062
$ws.Range("H132").Value = 250756.61
$ws.Range("I132").Value = 289977.78
$ws.Range("J132").Value = 15429.571
$ws.Range("K132").Value = 869933.3400000001
$ws.Range("L132").Value = 46288.713
$ws.Range("M132").Value = -867403.3400000001
$ws.Range("N132").Value = -51348.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14922.463
$ws.Range("I32").Value = 2194.5264
$ws.Range("K32").Value = 2194.5264
$ws.Range("M32").Value = -1907.5264

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 25379.8
$ws.Range("J100").Value = 25379.8
$ws.Range("L100").Value = 25379.8
$ws.Range("N100").Value = -27543.8
$ws.Range("H106").Value = 18563.143
$ws.Range("J106").Value = 18563.143
$ws.Range("L106").Value = 18563.143
$ws.Range("N106").Value = -21087.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 25666.666
$ws.Range("J88").Value = 25666.666
$ws.Range("L88").Value = 25666.666
$ws.Range("N88").Value = -26478.666
$ws.Range("H91").Value = 25666.666
$ws.Range("J91").Value = 25666.666
$ws.Range("L91").Value = 25666.666
$ws.Range("N91").Value = -28474.666
$ws.Range("H99").Value = 5690351.5
$ws.Range("I99").Value = 6954541
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 6954541
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -6953043
$ws.Range("N99").Value = -4496
$ws.Range("H100").Value = 56790
$ws.Range("J100").Value = 56790
$ws.Range("L100").Value = 56790
$ws.Range("N100").Value = -58954
$ws.Range("H112").Value = 24213.334
$ws.Range("J112").Value = 24213.334
$ws.Range("L112").Value = 24213.334
$ws.Range("N112").Value = -27167.334
$ws.Range("H119").Value = 49680.5
$ws.Range("J119").Value = 49680.5
$ws.Range("L119").Value = 49680.5
$ws.Range("N119").Value = -59356.5
$ws.Range("H126").Value = 5690351.5
$ws.Range("I126").Value = 6954541
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 20863623
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -20861153
$ws.Range("N126").Value = -9440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1319.3594
$ws.Range("I131").Value = 588
$ws.Range("J131").Value = 1381.339
$ws.Range("K131").Value = 1764
$ws.Range("L131").Value = 4144.017
$ws.Range("M131").Value = 3276
$ws.Range("N131").Value = -14224.017

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5532.8945
$ws.Range("I70").Value = 5432.593
$ws.Range("J70").Value = 5779.091
$ws.Range("K70").Value = 5432.593
$ws.Range("L70").Value = 5779.091
$ws.Range("M70").Value = -5162.593
$ws.Range("N70").Value = -6319.091
$ws.Range("H73").Value = 5532.8945
$ws.Range("I73").Value = 5432.593
$ws.Range("J73").Value = 5779.091
$ws.Range("K73").Value = 5432.593
$ws.Range("L73").Value = 5779.091
$ws.Range("M73").Value = -4496.593
$ws.Range("N73").Value = -7651.091
$ws.Range("H80").Value = 2573.3333
$ws.Range("I80").Value = 2345.4546
$ws.Range("J80").Value = 3200
$ws.Range("K80").Value = 2345.4546
$ws.Range("L80").Value = 3200
$ws.Range("M80").Value = -1347.4546
$ws.Range("N80").Value = -5196
$ws.Range("H83").Value = 2573.3333
$ws.Range("I83").Value = 2345.4546
$ws.Range("J83").Value = 3200
$ws.Range("K83").Value = 11727.273
$ws.Range("L83").Value = 16000
$ws.Range("M83").Value = -6735.273000000001
$ws.Range("N83").Value = -25984
$ws.Range("H122").Value = 301491
$ws.Range("I122").Value = 412351.88
$ws.Range("J122").Value = 2166.6
$ws.Range("K122").Value = 1237055.64
$ws.Range("L122").Value = 6499.799999999999
$ws.Range("M122").Value = -1234605.64
$ws.Range("N122").Value = -11399.8
$ws.Range("H126").Value = 2199.2307
$ws.Range("I126").Value = 1699
$ws.Range("K126").Value = 5097
$ws.Range("M126").Value = -2627

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2652.56
$ws.Range("I40").Value = 1629.25
$ws.Range("J40").Value = 3134.1177
$ws.Range("K40").Value = 1629.25
$ws.Range("L40").Value = 3134.1177
$ws.Range("M40").Value = -1493.25
$ws.Range("N40").Value = -3406.1177
$ws.Range("H94").Value = 7500
$ws.Range("J94").Value = 7500
$ws.Range("L94").Value = 7500
$ws.Range("N94").Value = -8852
$ws.Range("H106").Value = 22037.416
$ws.Range("J106").Value = 22037.416
$ws.Range("L106").Value = 22037.416
$ws.Range("N106").Value = -24561.416
$ws.Range("H122").Value = 3216.4614
$ws.Range("I122").Value = 2479
$ws.Range("J122").Value = 3757.2666
$ws.Range("K122").Value = 7437
$ws.Range("L122").Value = 11271.7998
$ws.Range("M122").Value = -4987
$ws.Range("N122").Value = -16171.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 168190.67
$ws.Range("J97").Value = 168190.67
$ws.Range("L97").Value = 168190.67
$ws.Range("N97").Value = -170172.67
$ws.Range("H112").Value = 25935.285
$ws.Range("J112").Value = 25935.285
$ws.Range("L112").Value = 25935.285
$ws.Range("N112").Value = -28889.285
$ws.Range("H122").Value = 78464.234
$ws.Range("I122").Value = 101104
$ws.Range("J122").Value = 2998.3333
$ws.Range("K122").Value = 303312
$ws.Range("L122").Value = 8994.999899999999
$ws.Range("M122").Value = -300862
$ws.Range("N122").Value = -13894.9999
$ws.Range("H125").Value = 28153.934
$ws.Range("J125").Value = 28153.934
$ws.Range("L125").Value = 28153.934
$ws.Range("N125").Value = -37993.934
$ws.Range("H126").Value = 59526.65
$ws.Range("I126").Value = 83772.75
$ws.Range("J126").Value = 1336
$ws.Range("K126").Value = 251318.25
$ws.Range("L126").Value = 4008
$ws.Range("M126").Value = -248848.25
$ws.Range("N126").Value = -8948
$ws.Range("H136").Value = 2236.92
$ws.Range("I136").Value = 1074.5625
$ws.Range("J136").Value = 4303.3335
$ws.Range("K136").Value = 3223.6875
$ws.Range("L136").Value = 12910.0005
$ws.Range("M136").Value = -673.6875
$ws.Range("N136").Value = -18010.0005
